$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196 (existing rows 196-211 shift down to 197-212)
$ws.Rows.Item(196).Insert()

# Populate the new row 196 with the new data point.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Categoria ID,
# G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo, L Precio maximo,
# M Precio promedio ponderado, N Unidad de comercializacion, O Origen,
# P Precio $/Kg, Q Kg o Unidades, R Clasificacion
$ws.Range("A196").Value = 8
$ws.Range("B196").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C196").Value = 'Coquimbo'
$ws.Range("D196").Value = 44746
$ws.Range("E196").Value = 4
$ws.Range("F196").Value = 100112037
$ws.Range("G196").Value = 'Cebollín'
$ws.Range("H196").Value = 'Sin especificar'
$ws.Range("I196").Value = 'Primera'
$ws.Range("J196").Value = 1000
$ws.Range("K196").Value = 1400
$ws.Range("L196").Value = 1600
$ws.Range("M196").Value = 1500
$ws.Range("N196").Value = '$/paquete 6 unidades'
$ws.Range("O196").Value = 'Provincia del Elquí'
$ws.Range("P196").Value = 250
$ws.Range("Q196").Value = 6
$ws.Range("R196").Value = 'Hortaliza'
